$d = $word.ActiveDocument

# Curly quotes used in "Change "understand" to "describe"".
$lq = [char]8220
$rq = [char]8221

$old = "If this class is not approved or does not fill then I will likely propose it again next year for a regular term (either in-person or on-line)."

# Paragraph count before the edit -- the new paragraphs land right after
# the paragraph that currently holds $old (the last paragraph in the doc).
$origCount = $d.Paragraphs.Count

# Split the final paragraph into five paragraphs by inserting paragraph
# marks (^p) plus the new text via Find/Replace -- this reliably creates
# real w:p boundaries (unlike InsertParagraphAfter at story end, which in
# this runtime clobbers the trailing paragraph).
$new = $old + "^p" + "^p" + `
       "Academic Council Meeting (8-Apr) Notes" + "^p" + `
       "Change " + $lq + "understand" + $rq + " to " + $lq + "describe" + $rq + "^p"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# After the split, the original last paragraph ("If this class...") keeps
# its index; the four new paragraphs follow it in order:
#   N+1 : blank separator paragraph
#   N+2 : "Academic Council Meeting (8-Apr) Notes" (bold heading)
#   N+3 : "Change "understand" to "describe"" (bulleted list item)
#   N+4 : paragraph that now holds the trailing bookmark
$blankPara = $d.Paragraphs.Item($origCount + 1)
$headingPara = $d.Paragraphs.Item($origCount + 2)
$changePara = $d.Paragraphs.Item($origCount + 3)
$bookmarkPara = $d.Paragraphs.Item($origCount + 4)

# Blank paragraph: strip the inherited bullet/list formatting entirely.
$blankPara.Range.ListFormat.RemoveNumbers()
$blankPara.Range.Style = "Normal"

# Heading paragraph: strip bullet/list formatting, make it bold.
$headingPara.Range.ListFormat.RemoveNumbers()
$headingPara.Range.Style = "Normal"
$headingPara.Range.Bold = 1

# "Change ... " paragraph keeps the inherited ListParagraph/bullet
# formatting (numId 1), matching the rest of the bulleted notes -- no
# further change needed there.

# Final paragraph (holds the _GoBack bookmark): strip bullet/list
# formatting so only the bookmark remains.
$bookmarkPara.Range.ListFormat.RemoveNumbers()
$bookmarkPara.Range.Style = "Normal"
